$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.725.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.119.72'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +10.22%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '255.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.667'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.50%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.62'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '60.70'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.374'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0740'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.58%  '
$ws.Range("E12").Value = '  +0.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.431.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +10.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.844'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.120.54'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +10.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.13'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.747.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0841'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.32'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '241.62'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("E25").Value = '  -7.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.99'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.57'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.64%  '
$ws.Range("E28").Value = '  +4.58%  '
$ws.Range("E29").Value = '  -9.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +67.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.123'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.52'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0960'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +11.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0600'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.35'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +16.19%  '
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.90'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.50%  '
$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.940'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.18'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.36'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.19'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0225'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.15'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.68%  '
$ws.Range("E44").Value = '  +13.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.357.77'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0843'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.12'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.308.83'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.33%  '
$ws.Range("E50").Value = '  -2.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.83'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.66%  '
